# Applies the "Additional iteration" edit to translator_testing_model.xlsx:
#  1. Inserts a new "Qualifier" sheet right after "TestEntityParameter",
#     carrying the same simple parameter/value header shape.
#  2. Splits the "predicate" column into "predicate_id" + "predicate_name"
#     on the three sheets that had a single "predicate" column:
#     TestAsset, AcceptanceTestAsset, TestEdgeData.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Qualifier" worksheet -------------------------------
$qualifier = $wb.Worksheets.Add($null, $wb.Worksheets.Item("TestEntityParameter"))
$qualifier.Name = "Qualifier"
$qualifier.Range("A1").Value = "parameter"
$qualifier.Range("B1").Value = "value"

# Match the look & feel (outline + margins) used by the rest of the workbook's sheets.
$qualifier.Outline.SummaryRow = 1
$qualifier.Outline.SummaryColumn = 1
$qualifier.PageSetup.LeftMargin = 54
$qualifier.PageSetup.RightMargin = 54
$qualifier.PageSetup.TopMargin = 72
$qualifier.PageSetup.BottomMargin = 72
$qualifier.PageSetup.HeaderMargin = 36
$qualifier.PageSetup.FooterMargin = 36

# --- 2. Split "predicate" -> "predicate_id" + "predicate_name" -------------

# TestAsset: predicate is column D (4)
$ws = $wb.Worksheets.Item("TestAsset")
$ws.Columns.Item(5).Insert()
$ws.Range("D1").Value = "predicate_id"
$ws.Range("E1").Value = "predicate_name"

# AcceptanceTestAsset: predicate is column N (14)
$ws = $wb.Worksheets.Item("AcceptanceTestAsset")
$ws.Columns.Item(15).Insert()
$ws.Range("N1").Value = "predicate_id"
$ws.Range("O1").Value = "predicate_name"

# TestEdgeData: predicate is column D (4)
$ws = $wb.Worksheets.Item("TestEdgeData")
$ws.Columns.Item(5).Insert()
$ws.Range("D1").Value = "predicate_id"
$ws.Range("E1").Value = "predicate_name"

# Restore the originally active sheet/selection.
$wb.Worksheets.Item("TestEntityParameter").Activate() | Out-Null
$wb.Worksheets.Item("TestEntityParameter").Range("A1").Select() | Out-Null
